$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported. Insert a new row at 254,
# pushing all the following rows (old 254-271) down by one (to 255-272),
# then populate the new row 254 with the new record's data. All of the
# descriptive columns (market, region, category, variety, quality, unit,
# origin, classification, etc.) are identical to the surrounding rows, so
# copy them from row 255 (the row that used to be row 254) into the new
# row 254, and then overwrite just the fields that actually changed.

$ws.Rows.Item(254).Insert()

for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(254, $col).Value = $ws.Cells.Item(255, $col).Value2
}

$ws.Cells.Item(254, 4).Value = 44714   # D254 Fecha
$ws.Cells.Item(254, 10).Value = 50     # J254 Volumen
$ws.Cells.Item(254, 11).Value = 22000  # K254 Precio minimo
$ws.Cells.Item(254, 12).Value = 22000  # L254 Precio maximo
$ws.Cells.Item(254, 13).Value = 22000  # M254 Precio promedio ponderado
$ws.Cells.Item(254, 16).Value = 2200   # P254 Precio $/Kg
